# Regional Availability Factor.xlsx - "Update to latest 4.0"
#
# - Switch the active sheet from "About" to "RAF-generation" and move the
#   selection there to B3.
# - Bump three RAF values up to 1 (solar thermal, biomass, natural gas peaker).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RAF-generation")

# Make RAF-generation the active sheet (drives tabSelected/activeTab).
$ws.Activate()

# Update the regional availability factors.
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 1
$ws.Range("B14").Value = 1

# Leave the selection on B3, matching the saved view state.
$ws.Range("B3").Select()
